# Overwrite figures with past ones with right size
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 51.7468893271431
$ws.Range("B3").Value = 54.2040008566919
$ws.Range("B4").Value = 55.1540879537179
$ws.Range("B5").Value = 55.8431322893747
